# Add 2022-Q4 data
# 1) Insert a new worksheet "2022-Q4" right after "总计" (so it ends up
#    before "2022-Q2"), populate it with the new quarterly fund-holding data.
# 2) Insert a new row into "总计" (summary) sheet for the 2022-Q4 totals,
#    shifting the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q4" sheet right after "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $total)
$q4.Name = "2022-Q4"

# Reuse the header formatting (bold + border) from an existing quarterly sheet
$headerStyle = $wb.Worksheets.Item("2022-Q1").Range("B1")
$headerStyle.Copy($q4.Range("B1:H1"))

$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Reuse the row-index column formatting (bold + border) too
$aStyle = $wb.Worksheets.Item("2022-Q1").Range("A2")
$aStyle.Copy($q4.Range("A2:A5"))

$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(5,1).Value = 3

# Force the numeric-looking columns to stay text (to match source formatting,
# e.g. fund codes keep leading zeros, percentages keep trailing zeros)
$q4.Range("B2:B5").NumberFormat = "@"
$q4.Range("D2:G5").NumberFormat = "@"

$q4.Cells.Item(2,2).Value = "001445"
$q4.Cells.Item(2,3).Value = "华安国企改革主题灵活配置混合A"
$q4.Cells.Item(2,4).Value = "15.25"
$q4.Cells.Item(2,5).Value = "85.79"
$q4.Cells.Item(2,6).Value = "3.19"
$q4.Cells.Item(2,7).Value = "0.4865"
$q4.Cells.Item(2,8).Value = 9

$q4.Cells.Item(3,2).Value = "163823"
$q4.Cells.Item(3,3).Value = "中银稳健策略灵活配置混合"
$q4.Cells.Item(3,4).Value = "2.02"
$q4.Cells.Item(3,5).Value = "31.34"
$q4.Cells.Item(3,6).Value = "1.05"
$q4.Cells.Item(3,7).Value = "0.0212"
$q4.Cells.Item(3,8).Value = 9

$q4.Cells.Item(4,2).Value = "006952"
$q4.Cells.Item(4,3).Value = "中银景元回报混合"
$q4.Cells.Item(4,4).Value = "1.29"
$q4.Cells.Item(4,5).Value = "24.13"
$q4.Cells.Item(4,6).Value = "1.12"
$q4.Cells.Item(4,7).Value = "0.0144"
$q4.Cells.Item(4,8).Value = 5

$q4.Cells.Item(5,2).Value = "016290"
$q4.Cells.Item(5,3).Value = "华安国企改革主题灵活配置混合C"
$q4.Cells.Item(5,4).Value = "0.16"
$q4.Cells.Item(5,5).Value = "85.79"
$q4.Cells.Item(5,6).Value = "3.19"
$q4.Cells.Item(5,7).Value = "0.0051"
$q4.Cells.Item(5,8).Value = 9

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q4 row into the "总计" summary sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows("2:2").Insert()

# The row-insert leaves stray formatting on the blank row; reset it back to
# the sheet's default (unstyled) look used by the other data rows.
$ws.Range("B2:D2").Style = "Normal"

# Give the new index cell (A2) the same bold/border style as the rows below.
$ws.Range("A3").Copy($ws.Range("A2"))

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "2022-Q4"
$ws.Cells.Item(2,3).Value = 4
$ws.Cells.Item(2,4).Value = 0.53

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(5,1).Value = 3
